$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (strikeout) values per row, replacing the prior "Strike#" derived values.
$kValues = @{
    "G2" = 2
    "G4" = 1
    "G5" = 4
    "G6" = 1
    "G7" = 1
    "G8" = 3
    "G9" = 2
    "G10" = 1
    "G11" = 0
    "G12" = 1
    "G13" = 1
    "G14" = 2
    "G15" = 3
    "G16" = 2
    "G17" = 1
    "G18" = 1
    "G19" = 3
    "G20" = 1
    "G21" = 0
    "G22" = 2
    "G23" = 2
    "G24" = 2
    "G25" = 3
    "G26" = 1
    "G27" = 3
    "G28" = 0
    "G29" = 1
    "G30" = 3
    "G31" = 5
    "G32" = 3
    "G33" = 3
    "G34" = 2
    "G35" = 2
    "G36" = 2
    "G37" = 1
    "G38" = 2
    "G39" = 0
    "G40" = 2
    "G41" = 0
    "G42" = 1
    "G43" = 1
    "G44" = 1
    "G45" = 1
    "G46" = 1
    "G47" = 2
    "G48" = 1
    "G49" = 4
    "G50" = 4
    "G51" = 1
    "G52" = 3
    "G53" = 3
    "G54" = 1
    "G55" = 3
    "G56" = 1
    "G57" = 2
    "G58" = 1
    "G59" = 0
}

foreach ($cellRef in $kValues.Keys) {
    $ws.Range($cellRef).Value = $kValues[$cellRef]
}
